# LTG L10 "Beyond the Circular Economy" - Part 1
# Commit: "updated LTG L10" -> the title-slide subtitle gains a trailing " I"
# (the lecture title becomes "Lecture 10: Beyond the Circular Economy I").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Slide 1 has two shapes: the big deck title ("CustomShape 1") and the
# lecture/author block ("CustomShape 2"). The lecture title is the first
# run of the first paragraph of CustomShape 2's text.
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$oldTitle = "Lecture 10: Beyond the Circular Economy"
$newTitle = "Lecture 10: Beyond the Circular Economy I"

$full = $tr.Text
$idx = $full.IndexOf($oldTitle)
if ($idx -ge 0) {
    $chars = $tr.Characters($idx + 1, $oldTitle.Length)
    $chars.Text = $newTitle
}
